# The single data row (row 2) in the "UserDetails" sheet previously held a
# stale/unused test-data record (Philip Bergstrom). This regenerates the
# row with the newest fake-data record produced by the test-data generator
# (Lexie Altenwerth), matching the columns:
#   A=firstName B=lastName C=userName D=password E=email F=mobileNumber

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Lexie"
$ws.Range("B2").Value = "Altenwerth"
$ws.Range("C2").Value = "LexieAltenwerth60400"
$ws.Range("D2").Value = "bnovld5v9"
$ws.Range("E2").Value = "bradford.schowalter@yahoo.com"
$ws.Range("F2").Value = "1-376-864-1461"
